# docs: update the pdf
#
# Removes the "modules__school__playStoreLink" and
# "modules__school__appStoreLink" columns from Table1 on the "Connectors"
# sheet, clears the corresponding header cells on the "Default Values"
# sheet (keeping their style), and switches the active sheet/selection
# over to "Default Values".

$wb = $excel.ActiveWorkbook

# --- Connectors sheet: drop the last two table columns -------------------
$connectors = $wb.Worksheets.Item("Connectors")
$table = $connectors.ListObjects.Item(1)

$table.ListColumns.Item(14).Delete()
$table.ListColumns.Item(13).Delete()

# Fully clear what used to be the M1/N1 header cells (value + style) so the
# sheet's used range shrinks back down to A1:L2.
$connectors.Range("M1:N1").Clear()

$connectors.Range("L4").Select() | Out-Null

# --- Default Values sheet: clear the now-unused header cells -------------
$defaults = $wb.Worksheets.Item("Default Values")
$defaults.Range("I1:J1").ClearContents()

$defaults.Range("D4").Select() | Out-Null

# The sheet that should end up active/selected is "Default Values".
$defaults.Activate() | Out-Null
